# Fixed update to excel issue
#
# The forecast was re-run a week later than the previous snapshot, so the
# "Forecast Comparison" sheet's 16-week window now starts 2025-02-02 (was
# 2025-01-26) and every forecast number is refreshed. The "Summary" sheet's
# derived statistics are updated to match.
#
# NOTE: every value touched below is stored as TEXT in the workbook (dates
# like "2025-02-02" and numbers like "157" are all plain strings, not real
# Excel dates/numbers). Excel's normal smart-input would silently convert
# such literals to a date serial or a number, so we prefix them with a
# leading apostrophe (exactly what typing `'2025-02-02` into a cell does)
# to force them to stay text, matching the original authoring.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------
# Columns: A=Week, B=Week_Start_Date, C=ASIN, D=MyForecast, E=Amazon Mean Forecast,
# F=Amazon P70 Forecast, G=Amazon P80 Forecast, H=Amazon P90 Forecast,
# I=Product Title, J=is_holiday_week  (A, C, I, J are unchanged by this update)

$forecastRows = @(
    @{ Row = 2;  Date = "2025-02-02"; D = 80; E = 124; F = 150; G = 178; H = 221 },
    @{ Row = 3;  Date = "2025-02-09"; D = 77; E = 104; F = 126; G = 150; H = 188 },
    @{ Row = 4;  Date = "2025-02-16"; D = 78; E = 111; F = 134; G = 161; H = 203 },
    @{ Row = 5;  Date = "2025-02-23"; D = 78; E = 111; F = 135; G = 162; H = 206 },
    @{ Row = 6;  Date = "2025-03-02"; D = 79; E = 112; F = 136; G = 165; H = 211 },
    @{ Row = 7;  Date = "2025-03-09"; D = 79; E = 112; F = 137; G = 167; H = 216 },
    @{ Row = 8;  Date = "2025-03-16"; D = 77; E = 108; F = 132; G = 162; H = 210 },
    @{ Row = 9;  Date = "2025-03-23"; D = 75; E = 108; F = 132; G = 163; H = 213 },
    @{ Row = 10; Date = "2025-03-30"; D = 74; E = 106; F = 130; G = 159; H = 206 },
    @{ Row = 11; Date = "2025-04-06"; D = 75; E = 103; F = 125; G = 154; H = 201 },
    @{ Row = 12; Date = "2025-04-13"; D = 75; E = 100; F = 122; G = 152; H = 200 },
    @{ Row = 13; Date = "2025-04-20"; D = 75; E = 100; F = 123; G = 152; H = 200 },
    @{ Row = 14; Date = "2025-04-27"; D = 74; E = 98;  F = 120; G = 148; H = 194 },
    @{ Row = 15; Date = "2025-05-04"; D = 73; E = 93;  F = 114; G = 143; H = 189 },
    @{ Row = 16; Date = "2025-05-11"; D = 72; E = 95;  F = 116; G = 145; H = 190 },
    @{ Row = 17; Date = "2025-05-18"; D = 71; E = 94;  F = 115; G = 144; H = 191 }
)

foreach ($entry in $forecastRows) {
    $r = $entry.Row
    $wsForecast.Cells.Item($r, 2).Value = "'" + $entry.Date
    $wsForecast.Cells.Item($r, 4).Value = $entry.D
    $wsForecast.Cells.Item($r, 5).Value = $entry.E
    $wsForecast.Cells.Item($r, 6).Value = $entry.F
    $wsForecast.Cells.Item($r, 7).Value = $entry.G
    $wsForecast.Cells.Item($r, 8).Value = $entry.H
}

# --- Summary sheet --------------------------------------------------------------
# All Value entries are text; numeric-looking / date-looking strings get a
# leading apostrophe so they remain text instead of becoming real numbers/dates.
$wsSummary.Cells.Item(2, 2).Value  = "2022-12-25 to 2025-01-26"
$wsSummary.Cells.Item(5, 2).Value  = "'157"
$wsSummary.Cells.Item(6, 2).Value  = "'137"
$wsSummary.Cells.Item(8, 2).Value  = "17109 units"
$wsSummary.Cells.Item(9, 2).Value  = "'1212"
$wsSummary.Cells.Item(10, 2).Value = "'622"
$wsSummary.Cells.Item(11, 2).Value = "'313"
$wsSummary.Cells.Item(12, 2).Value = "'80"
$wsSummary.Cells.Item(14, 2).Value = "'71"
$wsSummary.Cells.Item(15, 2).Value = "'2025-05-18"
